$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.920.34"
$ws.Range("E2").Value = "  -1.22%  "

$ws.Range("D3").Value = "1.914.05"
$ws.Range("E3").Value = "  +1.04%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'319.89"
$ws.Range("E5").Value = "  -1.00%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "'0.5033"
$ws.Range("E7").Value = "  -2.67%  "

$ws.Range("D8").Value = "'0.4037"
$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "'0.08317"
$ws.Range("E9").Value = "  -1.43%  "

$ws.Range("D10").Value = "'42.38"
$ws.Range("E10").Value = "  -0.87%  "

$ws.Range("D11").Value = "'1.101"
$ws.Range("E11").Value = "  -1.39%  "

$ws.Range("D12").Value = "'23.78"
$ws.Range("E12").Value = "  +1.28%  "

$ws.Range("D13").Value = "1.915.90"
$ws.Range("E13").Value = "  +1.25%  "

$ws.Range("D14").Value = "'6.392"
$ws.Range("E14").Value = "  -1.15%  "

$ws.Range("D15").Value = "'7.210"
$ws.Range("E15").Value = "  -1.58%  "

$ws.Range("D16").Value = "'1.005"
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D17").Value = "'91.99"
$ws.Range("E17").Value = "  -2.69%  "

$ws.Range("D18").Value = "'0.00001097"
$ws.Range("E18").Value = "  -1.07%  "

$ws.Range("D19").Value = "'0.06494"
$ws.Range("E19").Value = "  -2.36%  "

$ws.Range("D20").Value = "'18.19"
$ws.Range("E20").Value = "  -0.15%  "

$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").Value = "'5.925"
$ws.Range("E22").Value = "  -0.44%  "

$ws.Range("D23").Value = "29.962.27"

$ws.Range("D24").Value = "'11.28"
$ws.Range("E24").Value = "  -0.19%  "

$ws.Range("D25").Value = "'2.190"
$ws.Range("E25").Value = "  -1.53%  "

$ws.Range("D26").Value = "2.143.35"
$ws.Range("E26").Value = "  +1.54%  "

$ws.Range("D27").Value = "'21.93"
$ws.Range("E27").Value = "  +1.48%  "

$ws.Range("D28").Value = "'161.86"
$ws.Range("E28").Value = "  -0.52%  "

$ws.Range("D29").Value = "'2.319"
$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("D30").Value = "'128.62"
$ws.Range("E30").Value = "  -0.75%  "

$ws.Range("D31").Value = "'1.140"
$ws.Range("E31").Value = "  +4.66%  "

$ws.Range("D32").Value = "'0.1035"
$ws.Range("E32").Value = "  -1.58%  "

$ws.Range("D33").Value = "'5.947"

$ws.Range("D34").Value = "'3.737"
$ws.Range("E34").Value = "  -0.30%  "

$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "'0.02439"
$ws.Range("E35").Value = "  -2.01%  "

$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'5.378"
$ws.Range("E36").Value = "  +0.82%  "

$ws.Range("D37").Value = "'0.06378"
$ws.Range("E37").Value = "  -2.38%  "

$ws.Range("D38").Value = "'0.2144"
$ws.Range("E38").Value = "  -2.70%  "

$ws.Range("D39").Value = "'0.6539"
$ws.Range("E39").Value = "  +0.59%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'1.189"
$ws.Range("E40").Value = "  -2.28%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'8.659"
$ws.Range("E41").Value = "  -1.52%  "

$ws.Range("D42").Value = "'11.32"
$ws.Range("E42").Value = "  -4.93%  "

$ws.Range("D43").Value = "'1.209"
$ws.Range("E43").Value = "  -1.79%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'13.50"
$ws.Range("E44").Value = "  +2.44%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'2.204"
$ws.Range("E45").Value = "  +6.87%  "

$ws.Range("D46").Value = "'0.6053"
$ws.Range("E46").Value = "  -0.78%  "

$ws.Range("D47").Value = "'3.615"
$ws.Range("E47").Value = "  -1.83%  "

$ws.Range("D48").Value = "'1.206"
$ws.Range("E48").Value = "  -2.54%  "

$ws.Range("D49").Value = "'121.36"
$ws.Range("E49").Value = "  -3.00%  "

$ws.Range("D50").Value = "'78.77"
$ws.Range("E50").Value = "  -0.54%  "

$ws.Range("D51").Value = "'1.129"
$ws.Range("E51").Value = "  -1.78%  "
